# Update Financials: insert two new quarterly columns (most-recent quarters)
# before column D, shifting the existing quarter columns right, then populate
# the two new columns with the newly reported figures, and correct a couple
# of restated prior-quarter figures that came in with the refreshed data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 2 blank columns before column D (pushes D:K -> F:M).
$ws.Range("D1:E1").EntireColumn.Insert()

# 2) Give the two new columns (D:E) the same number formatting / font as the
#    columns that were just shifted right (F:M), so the new quarters look
#    like the rest of the table (date header row uses style 2, data rows
#    use style 3 - PasteSpecial formats grabs both from the range below).
$ws.Range("F7:M102").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Populate the two new columns (D = most recent quarter, E = prior
#    quarter) for every row that carries data in this report. $null left
#    where the source row has no figure (section headers / spacer rows),
#    matching the blank-but-styled cells elsewhere in the sheet.
$newQuarterData = @{
    7 = @(43465, 43373);
    8 = @(11543000, 11412000);
    9 = @(11122000, 10494000);
    10 = @(421000, 918000);
    11 = @($null, $null);
    12 = @("NA", "NA");
    13 = @(0, 0);
    14 = @(0, "NA");
    15 = @(0, 0);
    16 = @($null, $null);
    17 = @(11490000, 10827000);
    18 = @(53000, 585000);
    19 = @($null, $null);
    20 = @(34000, -35000);
    21 = @(246000, 709000);
    22 = @(70000, 98000);
    23 = @(17000, 452000);
    24 = @(99000, 85000);
    25 = @(0, 0);
    26 = @(-82000, 367000);
    27 = @(-98000, 350000);
    28 = @(0, 0);
    29 = @(24000, 7000);
    30 = @(0, 0);
    31 = @(0, 0);
    32 = @(-34000, 35000);
    33 = @(-74000, 357000);
    34 = @(0, 0);
    35 = @(-74000, 357000);
    38 = @(43465, 43373);
    39 = @($null, $null);
    40 = @($null, $null);
    41 = @(389000, 267000);
    42 = @(510000, 629000);
    43 = @(2367000, 2377000);
    44 = @(5871000, 7183000);
    45 = @(1931000, 2655000);
    46 = @(11068000, 13111000);
    47 = @(1120000, 1133000);
    48 = @(5201000, 5164000);
    49 = @(1424000, 1434000);
    50 = @(0, 0);
    51 = @(0, 0);
    52 = @(612000, 604000);
    53 = @(0, 0);
    54 = @(19425000, 21446000);
    55 = @($null, $null);
    56 = @($null, $null);
    57 = @(3501000, 3274000);
    58 = @(1169000, 2388000);
    59 = @(2502000, 2639000);
    60 = @(7172000, 8301000);
    61 = @(4203000, 4912000);
    62 = @(1248000, 1308000);
    63 = @(0, 0);
    64 = @(0, 0);
    65 = @(0, 0);
    66 = @(13252000, 15166000);
    67 = @($null, $null);
    68 = @(0, 0);
    69 = @(0, 0);
    70 = @(690000, 690000);
    71 = @(0, 0);
    72 = @(8059000, 8203000);
    73 = @(0, 0);
    74 = @(0, 0);
    75 = @(0, 0);
    76 = @(5483000, 5590000);
    77 = @(0, 0);
    80 = @(43465, 43373);
    81 = @(-74000, 357000);
    82 = @($null, $null);
    83 = @(159000, 159000);
    84 = @(0, 0);
    85 = @(0, 0);
    86 = @(0, 0);
    87 = @(0, 0);
    88 = @(0, 0);
    89 = @(2021000, 53000);
    90 = @($null, $null);
    91 = @(-175000, -98000);
    92 = @(0, 0);
    93 = @(0, 0);
    94 = @(186000, 652000);
    95 = @($null, $null);
    96 = @(-46000, -78000);
    97 = @(0, 0);
    98 = @(0, 0);
    99 = @(0, 0);
    100 = @(-2078000, -709000);
    101 = @(-7000, 50000);
    102 = @(122000, 46000)
}

foreach ($row in $newQuarterData.Keys) {
    $pair = $newQuarterData[$row]
    $ws.Cells.Item($row, 4).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
}

# 4) A couple of quarters further back in the "Dividends Paid" and "Other
#    Cashflows from Investing Activities" rows were restated in this data
#    refresh (not just shifted) - correct those explicitly.
$ws.Cells.Item(91, 9).Value = -143000   # I91
$ws.Cells.Item(91, 10).Value = -160000  # J91
$ws.Cells.Item(94, 8).Value = 1322000   # H94
$ws.Cells.Item(94, 9).Value = 502000    # I94
